$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.365.05"
$ws.Range("E2").Value = "  -4.61%  "
$ws.Range("D3").Value = "3.266.25"
$ws.Range("E3").Value = "  -7.05%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'594.80"
$ws.Range("E5").Value = "  -4.73%  "
$ws.Range("D6").Value = "'151.31"
$ws.Range("E6").Value = "  -12.04%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.256.88"
$ws.Range("E8").Value = "  -7.19%  "
$ws.Range("D9").Value = "'0.541"
$ws.Range("E9").Value = "  -11.16%  "
$ws.Range("D10").Value = "'0.171"
$ws.Range("E10").Value = "  -14.16%  "
$ws.Range("D11").Value = "'6.62"
$ws.Range("E11").Value = "  -7.82%  "
$ws.Range("D12").Value = "'0.510"
$ws.Range("E12").Value = "  -12.89%  "
$ws.Range("D13").Value = "'38.11"
$ws.Range("E13").Value = "  -17.50%  "
$ws.Range("D14").Value = "'0.0000244"
$ws.Range("E14").Value = "  -11.63%  "
$ws.Range("D15").Value = "3.788.15"
$ws.Range("E15").Value = "  -7.16%  "
$ws.Range("D16").Value = "67.354.20"
$ws.Range("E16").Value = "  -4.75%  "
$ws.Range("D17").Value = "3.267.97"
$ws.Range("E17").Value = "  -6.95%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'535.78"
$ws.Range("E18").Value = "  -11.70%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.114"
$ws.Range("E19").Value = "  -6.46%  "
$ws.Range("D20").Value = "'7.21"
$ws.Range("E20").Value = "  -14.23%  "
$ws.Range("D21").Value = "'15.06"
$ws.Range("E21").Value = "  -15.02%  "
$ws.Range("D22").Value = "'0.761"
$ws.Range("E22").Value = "  -13.45%  "
$ws.Range("D23").Value = "'7.88"
$ws.Range("E23").Value = "  -13.15%  "
$ws.Range("D24").Value = "'85.46"
$ws.Range("E24").Value = "  -11.96%  "
$ws.Range("D25").Value = "'13.56"
$ws.Range("E25").Value = "  -12.56%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E27").Value = "  -12.77%  "
$ws.Range("D28").Value = "'29.34"
$ws.Range("E28").Value = "  -12.33%  "
$ws.Range("D29").Value = "'8.01"
$ws.Range("E29").Value = "  -11.20%  "
$ws.Range("D30").Value = "'2.13"
$ws.Range("E30").Value = "  -16.68%  "
$ws.Range("E31").Value = "  -11.69%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  -11.11%  "
$ws.Range("D33").Value = "'542.45"
$ws.Range("E33").Value = "  -12.68%  "
$ws.Range("D34").Value = "'6.62"
$ws.Range("E34").Value = "  -17.87%  "
$ws.Range("D35").Value = "'5.70"
$ws.Range("E35").Value = "  -15.99%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'0.0454"
$ws.Range("E37").Value = "  -7.30%  "
$ws.Range("D38").Value = "'53.21"
$ws.Range("E38").Value = "  -5.95%  "
$ws.Range("D39").Value = "'0.0857"
$ws.Range("E39").Value = "  -13.85%  "
$ws.Range("E40").Value = "  -10.02%  "
$ws.Range("D41").Value = "'9.06"
$ws.Range("E41").Value = "  -16.14%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.72"
$ws.Range("E42").Value = "  -20.36%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.931.65"
$ws.Range("E43").Value = "  -12.17%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0590"
$ws.Range("E44").Value = "  -18.25%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.262"
$ws.Range("E45").Value = "  -15.63%  "
$ws.Range("D46").Value = "'26.83"
$ws.Range("E46").Value = "  -15.68%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.14"
$ws.Range("E47").Value = "  -14.38%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'127.07"
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("D50").Value = "'2.32"
$ws.Range("E50").Value = "  -21.09%  "
$ws.Range("E51").Value = "  -12.83%  "
